# "fix error on questions"
# Fill in missing/incorrect question data (id numbers, video src ids,
# startTime and answerKeyword corrections) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- answerKeyword (column E) corrections / additions ---
$ws.Range("E4").Value  = "길크 길티"
$ws.Range("E16").Value = "아이마스 신데마스 아이돌마스터"
$ws.Range("E20").Value = "메이도라 메이드래곤 코바야시"
$ws.Range("E24").Value = "카구야 카구야님"
$ws.Range("E26").Value = "바이올렛"
$ws.Range("E35").Value = "에망센 에로망가"
$ws.Range("E47").Value = "신만세"
$ws.Range("E50").Value = "엔비"
$ws.Range("E51").Value = "기교소녀"

# --- answer (column D) text fix ---
$ws.Range("D48").Value = "럭키스타"

# --- row 10: wrong video src id / startTime ---
$ws.Range("B10").Value = "CFM_zypYFHM"
$ws.Range("C10").Value = 50

# --- rows 39-52: fill in missing id (A), src (B) and startTime (C) ---
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = "NUGqGRcNb38"
$ws.Range("C39").Value = 0

$ws.Range("A40").Value = 39
$ws.Range("B40").Value = "49X8c3Q5Ezs"
$ws.Range("C40").Value = 0

$ws.Range("A41").Value = 40

$ws.Range("A42").Value = 41

$ws.Range("A43").Value = 42

$ws.Range("A44").Value = 43

$ws.Range("A45").Value = 44
$ws.Range("B45").Value = "wT5uF0kAW9Y"
$ws.Range("C45").Value = 0

$ws.Range("A46").Value = 45

$ws.Range("A47").Value = 46
$ws.Range("B47").Value = "v9iz9glzbMs"
$ws.Range("C47").Value = 200

$ws.Range("A48").Value = 47

$ws.Range("A49").Value = 48

$ws.Range("A50").Value = 49

$ws.Range("A51").Value = 50

$ws.Range("A52").Value = 51

# --- restore selection to B12 (matches author's last-saved cursor) ---
$ws.Range("B12").Select()
